$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.17%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.27%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.042"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.40%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08042"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.48%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.898"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.28%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.778"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.26%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9211"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.64%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1279"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.47%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1910"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.42%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09068"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.20%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03442"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.80%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09845"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.31%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001403"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.83%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006245"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "7.95%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.797"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.73%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.148"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.98%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "12.56%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3419"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.04%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.11%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.187"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.92%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2303"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-11.56%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04435"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.99%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.98%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004611"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.02%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001252"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-3.79%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "42.01%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01948"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.67%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05390"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "9.93%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007613"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.23%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01012"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.15%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.51%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002174"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.40%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009640"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.93%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006138"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.47%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.60%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001660"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.30%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.05%"
